# ENACOM_GH.xlsx — "Add files via upload" edit
#
# 1) Update the B12 input value (Hoja1).
# 2) Let D12's "C+B" formula recompute against the new B12 (it already
#    reads =C12+B12 via the shared formula inherited from D3, so simply
#    re-asserting the formula keeps it semantically identical while
#    picking up the new total).
# 3) Move the active selection on Hoja1 from C20 to C6.
#
# (The workbookView window-size metadata in xl/workbook.xml reflects the
# pixel size of the Excel application window at save time; it is not
# exposed anywhere on the Application/Window COM object model, so it
# can't be driven from automation code here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# New value for B12
$ws.Range("B12").Value = 5422117878

# Re-assert D12's total so the dependent formula result is refreshed.
$ws.Range("D12").Formula = "=C12+B12"

# Move the selection to C6
$ws.Range("C6").Select()
